$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 30; existing rows 30-41 shift down to 32-43.
$ws.Rows("30:31").Insert()

# Populate the two newly inserted rows (30 and 31) with new weekly data.
# Columns: A Mercado ID, B Mercado, C Región, D Fecha, E Codreg, F Tipo,
#          G Producto ID, H Producto, I Categoría ID, J Categoría, K Variedad,
#          L Calidad, M Volumen, N Precio mínimo, O Precio máximo,
#          P Precio promedio ponderado, Q Unidad de comercialización,
#          R Origen, S Precio $/Kg, T Kg / unidad

$newRows = @(
    @{ Row = 30; Fecha = 45001; Calidad = 'Especial'; Volumen = 150; PMin = 13000; PMax = 13000; PProm = 13000; Unidad = '$/caja 18 kilos'; Origen = 'Provincia de Melipilla'; PKg = 722; KgUnidad = 18 },
    @{ Row = 31; Fecha = 45001; Calidad = 'Primera';  Volumen = 100; PMin = 11000; PMax = 11000; PProm = 11000; Unidad = '$/caja 18 kilos'; Origen = 'Provincia de Melipilla'; PKg = 611; KgUnidad = 18 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 5
    $ws.Cells.Item($row, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($row, 3).Value = "Maule"
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 5).Value = 7
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100107
    $ws.Cells.Item($row, 8).Value = "Otros"
    $ws.Cells.Item($row, 9).Value = 100107011
    $ws.Cells.Item($row, 10).Value = "Tuna"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
